# Apply updated crypto price/volume figures per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.866.88"
$ws.Range("E2").Value = "  -1.05%  "

$ws.Range("D3").Value = "1.639.61"
$ws.Range("E3").Value = "  -0.79%  "

$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = "  +0.27%  "

$ws.Range("D5").Value = "'215.85"
$ws.Range("E5").Value = "  -0.07%  "

$ws.Range("D6").Value = "'0.5032"
$ws.Range("E6").Value = "  -1.47%  "

$ws.Range("D7").Value = "'1.004"
$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "'0.2570"
$ws.Range("E8").Value = "  -0.50%  "

$ws.Range("D9").Value = "'0.06391"
$ws.Range("E9").Value = "  -0.57%  "

$ws.Range("D10").Value = "'19.72"
$ws.Range("E10").Value = "  -1.22%  "

$ws.Range("D11").Value = "'0.07728"
$ws.Range("E11").Value = "  -0.81%  "

$ws.Range("D12").Value = "1.653.06"
$ws.Range("E12").Value = "  +0.00%  "

$ws.Range("D13").Value = "'4.275"
$ws.Range("E13").Value = "  -0.08%  "

$ws.Range("D14").Value = "1.862.33"
$ws.Range("E14").Value = "  -0.96%  "

$ws.Range("D15").Value = "'0.5458"
$ws.Range("E15").Value = "  -1.08%  "

$ws.Range("D16").Value = "0.0₅7928"
$ws.Range("E16").Value = "  -1.29%  "

$ws.Range("D17").Value = "'63.88"
$ws.Range("E17").Value = "  -0.20%  "

$ws.Range("D18").Value = "25.883.48"
$ws.Range("E18").Value = "  -1.04%  "

$ws.Range("D19").Value = "'1.005"
$ws.Range("E19").Value = "  -0.05%  "

$ws.Range("D20").Value = "'201.35"
$ws.Range("E20").Value = "  -4.36%  "

$ws.Range("D21").Value = "'4.360"
$ws.Range("E21").Value = "  -1.16%  "

$ws.Range("D22").Value = "'9.932"
$ws.Range("E22").Value = "  -1.29%  "

$ws.Range("D23").Value = "'5.972"
$ws.Range("E23").Value = "  -1.07%  "

$ws.Range("D24").Value = "'1.005"
$ws.Range("E24").Value = "  +0.05%  "

$ws.Range("D25").Value = "'1.926"
$ws.Range("E25").Value = "  +10.86%  "

$ws.Range("D26").Value = "'141.71"
$ws.Range("E26").Value = "  -1.29%  "

$ws.Range("D27").Value = "'0.1136"
$ws.Range("E27").Value = "  -3.79%  "

$ws.Range("D28").Value = "'15.70"
$ws.Range("E28").Value = "  -0.75%  "

$ws.Range("D29").Value = "'6.721"
$ws.Range("E29").Value = "  -3.70%  "

$ws.Range("D30").Value = "'1.245"
$ws.Range("E30").Value = "  +0.02%  "

$ws.Range("D31").Value = "'0.04998"
$ws.Range("E31").Value = "  -2.02%  "

$ws.Range("D32").Value = "'3.278"
$ws.Range("E32").Value = "  -1.90%  "

$ws.Range("D33").Value = "'3.199"
$ws.Range("E33").Value = "  -0.68%  "

$ws.Range("D34").Value = "'1.542"
$ws.Range("E34").Value = "  -1.53%  "

$ws.Range("D35").Value = "'2.376"
$ws.Range("E35").Value = "  +0.63%  "

$ws.Range("D36").Value = "1.171.92"
$ws.Range("E36").Value = "  +0.31%  "

$ws.Range("D37").Value = "'2.632"
$ws.Range("E37").Value = "  -4.18%  "

$ws.Range("D38").Value = "'0.8925"
$ws.Range("E38").Value = "  -3.57%  "

$ws.Range("D39").Value = "'0.5587"
$ws.Range("E39").Value = "  -1.77%  "

$ws.Range("D40").Value = "'0.01561"
$ws.Range("E40").Value = "  -1.68%  "

$ws.Range("D41").Value = "'1.005"
$ws.Range("E41").Value = "  +0.08%  "

$ws.Range("D42").Value = "'5.697"
$ws.Range("E42").Value = "  +0.78%  "

$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'99.85"
$ws.Range("E43").Value = "  -0.52%  "

$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'0.8066"
$ws.Range("E44").Value = "  -2.96%  "

$ws.Range("D45").Value = "1.770.14"
$ws.Range("E45").Value = "  -1.14%  "

$ws.Range("E46").Value = "  -0.26%  "

$ws.Range("D47").Value = "'0.4529"
$ws.Range("E47").Value = "  -0.45%  "

$ws.Range("D48").Value = "'1.002"
$ws.Range("E48").Value = "  -0.43%  "

$ws.Range("D49").Value = "'54.95"
$ws.Range("E49").Value = "  -1.23%  "

$ws.Range("D50").Value = "'0.05073"
$ws.Range("E50").Value = "  +0.46%  "

$ws.Range("D51").Value = "'1.003"
$ws.Range("E51").Value = "  -0.31%  "
